# Regenerate save_data column G ("K") with updated strikeout values,
# replacing the previous Strike# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 6
    3  = 6
    4  = 4
    5  = 7
    6  = 14
    7  = 6
    8  = 12
    9  = 11
    10 = 7
    11 = 7
    12 = 1
    13 = 10
    14 = 5
    15 = 3
    16 = 6
    17 = 4
    18 = 3
    19 = 3
    20 = 7
    21 = 6
    22 = 6
    23 = 6
    24 = 5
    25 = 9
    26 = 3
    27 = 6
    28 = 6
    29 = 4
    30 = 4
    31 = 7
    32 = 4
    33 = 10
    34 = 9
    35 = 3
    36 = 5
    37 = 9
    38 = 4
    39 = 6
    40 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
